# Generate Report for Archive
#
# - Status text "Ready for handoff" -> "In Translation" on all three sheets
# - Narrower Status-related columns (Overview!E:F, zh-cn!C, de-de!C), matching
#   the new (shorter) status text, ~13.41 characters wide.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E1:F1").ColumnWidth = 12.5

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C1").ColumnWidth = 12.5

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C1").ColumnWidth = 12.5
